$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.801.61"
$ws.Range("E2").Value = "'  -0.29%  "
$ws.Range("D3").Value = "'3.410.85"
$ws.Range("E3").Value = "'  -0.10%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.42%  "
$ws.Range("D5").Value = "'412.60"
$ws.Range("E5").Value = "'  +0.92%  "
$ws.Range("D6").Value = "'129.29"
$ws.Range("E6").Value = "'  +0.30%  "
$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "'  -2.34%  "
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("D9").Value = "'0.726"
$ws.Range("E9").Value = "'  -1.09%  "
$ws.Range("D10").Value = "'0.140"
$ws.Range("E10").Value = "'  -1.17%  "
$ws.Range("D11").Value = "'42.68"
$ws.Range("E11").Value = "'  +0.71%  "
$ws.Range("D12").Value = "'0.0000220"
$ws.Range("E12").Value = "'  +2.80%  "
$ws.Range("D13").Value = "'9.16"
$ws.Range("E13").Value = "'  +1.95%  "
$ws.Range("D14").Value = "'3.949.87"
$ws.Range("E14").Value = "'  -0.34%  "
$ws.Range("E15").Value = "'  -0.21%  "
$ws.Range("D16").Value = "'20.48"
$ws.Range("E16").Value = "'  -1.62%  "
$ws.Range("D17").Value = "'3.423.82"
$ws.Range("E17").Value = "'  -0.02%  "
$ws.Range("D18").Value = "'12.72"
$ws.Range("E18").Value = "'  +5.12%  "
$ws.Range("D20").Value = "'61.838.55"
$ws.Range("E20").Value = "'  -0.27%  "
$ws.Range("D21").Value = "'478.09"
$ws.Range("E21").Value = "'  +8.07%  "
$ws.Range("D22").Value = "'90.87"
$ws.Range("E22").Value = "'  -0.13%  "
$ws.Range("D23").Value = "'3.27"
$ws.Range("E23").Value = "'  +3.47%  "
$ws.Range("D24").Value = "'13.07"
$ws.Range("E24").Value = "'  +0.49%  "
$ws.Range("D25").Value = "'3.31"
$ws.Range("E25").Value = "'  +2.12%  "
$ws.Range("D26").Value = "'9.72"
$ws.Range("E26").Value = "'  +10.47%  "
$ws.Range("D27").Value = "'33.14"
$ws.Range("E27").Value = "'  -1.61%  "
$ws.Range("E28").Value = "'  +0.27%  "
$ws.Range("D29").Value = "'7.69"
$ws.Range("E29").Value = "'  +1.46%  "
$ws.Range("D30").Value = "'11.84"
$ws.Range("E30").Value = "'  -0.82%  "
$ws.Range("E31").Value = "'  -3.08%  "
$ws.Range("E32").Value = "'  -1.98%  "
$ws.Range("E33").Value = "'  -2.98%  "
$ws.Range("D34").Value = "'40.90"
$ws.Range("E34").Value = "'  -3.61%  "
$ws.Range("D36").Value = "'58.14"
$ws.Range("E36").Value = "'  +7.91%  "
$ws.Range("D37").Value = "'0.0486"
$ws.Range("E37").Value = "'  -2.73%  "
$ws.Range("E38").Value = "'  +0.05%  "
$ws.Range("D39").Value = "'3.02"
$ws.Range("E39").Value = "'  +3.95%  "
$ws.Range("D40").Value = "'0.324"
$ws.Range("E40").Value = "'  +3.33%  "
$ws.Range("D41").Value = "'148.73"
$ws.Range("E41").Value = "'  +5.40%  "
$ws.Range("E42").Value = "'  -0.46%  "
$ws.Range("D43").Value = "'3.32"
$ws.Range("E43").Value = "'  -1.00%  "
$ws.Range("E44").Value = "'  +5.86%  "
$ws.Range("E45").Value = "'  +7.45%  "
$ws.Range("D46").Value = "'4.22"
$ws.Range("E46").Value = "'  +3.16%  "
$ws.Range("D47").Value = "'2.35"
$ws.Range("E47").Value = "'  +19.16%  "
$ws.Range("D48").Value = "'16.38"
$ws.Range("E48").Value = "'  -1.03%  "
$ws.Range("D49").Value = "'0.0₃0536"
$ws.Range("E49").Value = "'  +22.50%  "
$ws.Range("D50").Value = "'22.21"
$ws.Range("E50").Value = "'  +0.43%  "
$ws.Range("D51").Value = "'113.59"
$ws.Range("E51").Value = "'  +9.95%  "
